$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).RowHeight = 15.15

$ws.Range("A6").Value = "jean"
$ws.Range("B6").Value = 67347422
$ws.Range("C6").Value = "SIL (Système Informatique et Logiciel)"
$ws.Range("D6").Value = "cotonou"

$ws.Rows.Item(7).OutlineLevel = 5
$ws.Rows.Item(7).Delete()

$ws.Range("F10").Select()
